# Auto-generated from the OOXML diff: update crypto price/volume data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.837.57'
$ws.Range('E2').Value = '  +0.79%  '
$ws.Range('D3').Value = '1.707.15'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('D4').Value = '''1.002'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''316.61'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '''1.001'
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '''0.3940'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '''0.4057'
$ws.Range('E8').Value = '  +0.36%  '
$ws.Range('E9').Value = '  -1.73%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('E11').Value = '  -2.36%  '
$ws.Range('D12').Value = '''0.08917'
$ws.Range('E12').Value = '  +1.23%  '
$ws.Range('D13').Value = '''7.318'
$ws.Range('E13').Value = '  +0.56%  '
$ws.Range('D14').Value = '''23.63'
$ws.Range('E14').Value = '  +0.82%  '
$ws.Range('D15').Value = '''8.045'
$ws.Range('E15').Value = '  +5.36%  '
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').Value = '1.704.05'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '''100.41'
$ws.Range('E18').Value = '  -0.94%  '
$ws.Range('D19').Value = '''0.07049'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').Value = '''19.75'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').Value = '''7.087'
$ws.Range('E21').Value = '  +2.44%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('D23').Value = '''14.58'
$ws.Range('E23').Value = '  +2.76%  '
$ws.Range('D24').Value = '24.815.40'
$ws.Range('E24').Value = '  +0.77%  '
$ws.Range('D25').Value = '''3.222'
$ws.Range('E25').Value = '  +7.80%  '
$ws.Range('E26').Value = '  +1.44%  '
$ws.Range('D27').Value = '''22.89'
$ws.Range('E27').Value = '  +2.08%  '
$ws.Range('D28').Value = '''161.98'
$ws.Range('E28').Value = '  +1.40%  '
$ws.Range('D29').Value = '''8.389'
$ws.Range('E29').Value = '  +9.99%  '
$ws.Range('D30').Value = '''137.31'
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('D31').Value = '''5.176'
$ws.Range('E31').Value = '  -1.46%  '
$ws.Range('D32').Value = '''0.08914'
$ws.Range('E32').Value = '  +3.81%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D33').Value = '''7.561'
$ws.Range('E33').Value = '  +2.06%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = '''1.087'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').Value = '''11.18'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('D36').Value = '''1.979'
$ws.Range('E36').Value = '  +1.47%  '
$ws.Range('D37').Value = '''0.2757'
$ws.Range('E37').Value = '  -0.40%  '
$ws.Range('D38').Value = '''14.52'
$ws.Range('E38').Value = '  -1.70%  '
$ws.Range('D39').Value = '''0.09238'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('D40').Value = '''0.02762'
$ws.Range('E40').Value = '  -1.59%  '
$ws.Range('E41').Value = '  -0.33%  '
$ws.Range('D42').Value = '''0.7731'
$ws.Range('E42').Value = '  -0.36%  '
$ws.Range('D43').Value = '''15.86'
$ws.Range('E43').Value = '  +1.04%  '
$ws.Range('D44').Value = '''0.7225'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').Value = '''2.577'
$ws.Range('E45').Value = '  +2.33%  '
$ws.Range('D46').Value = '''4.211'
$ws.Range('E46').Value = '  -0.26%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '''140.70'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').Value = '''1.324'
$ws.Range('E49').Value = '  -5.83%  '
$ws.Range('D50').Value = '''91.22'
$ws.Range('E50').Value = '  +3.26%  '
$ws.Range('D51').Value = '''0.08015'
$ws.Range('E51').Value = '  -0.41%  '
